$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style into the new H column header before other edits
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Delete rows 3,4,5 (old extra data rows that are being removed)
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()

# Update header row
$ws.Range("B1").Value = "Subarea"
$ws.Range("C1").Value = "Location"
$ws.Range("D1").Value = "Approved CCTV Vendor"
$ws.Range("E1").Value = "Warranty"
$ws.Range("F1").Value = "OC Surveyor"
$ws.Range("G1").Value = "Reviewer.1"
$ws.Range("H1").Value = "Notes.1"

# Update data row 2
$ws.Range("B2").Value = "gh"
$ws.Range("C2").Value = "gh"
$ws.Range("D2").Value = "GPH"
$ws.Range("E2").Value = "Accepted"
$ws.Range("F2").Value = "COSTELLO_C"
$ws.Range("G2").Value = "T. Martin"
$ws.Range("H2").Value = "hg"

# Column widths (ColumnWidth is in "characters"; COM quantizes to whole
# pixels, so feed it the pixel-rounded character width that re-serializes
# to the target OOXML width)
$ws.Columns.Item(8).ColumnWidth = 13.5
$ws.Columns.Item(9).ColumnWidth = 26.5

# Narrow the "no blanks" conditional formatting range from the old B2:I5
# down to the new single data row, B2:L2
$fcs = $ws.Cells.FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    if ($fc.Type -eq 10) {
        $fc.ModifyAppliesToRange($ws.Range("B2:L2"))
    }
}

$wb.Save()
